$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.403.66"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").Value = "1.829.46"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.99"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4479"
$ws.Range("E7").Value = "  +5.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3778"
$ws.Range("E8").Value = "  +3.29%  "

$ws.Range("E9").Value = "  +4.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8995"
$ws.Range("E10").Value = "  +7.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.05"
$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("D12").Value = "1.822.61"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.775"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.46"
$ws.Range("E14").Value = "  +5.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.416"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07111"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008827"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9991"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").Value = "  +3.01%  "

$ws.Range("D21").Value = "27.420.56"
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.295"
$ws.Range("E22").Value = "  +3.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").Value = "2.060.07"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.006"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.492"
$ws.Range("E26").Value = "  +11.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.64"
$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.68"
$ws.Range("E28").Value = "  +2.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.396"
$ws.Range("E29").Value = "  +3.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.21"
$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08857"
$ws.Range("E31").Value = "  +1.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7829"
$ws.Range("E32").Value = "  +7.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.208"
$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.581"
$ws.Range("E34").Value = "  +4.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.886"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9988"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.112"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("E38").Value = "  +2.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05344"
$ws.Range("E39").Value = "  +2.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.394"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5350"
$ws.Range("E41").Value = "  +4.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1732"
$ws.Range("E42").Value = "  +3.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.867"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.284"
$ws.Range("E44").Value = "  +17.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.834"
$ws.Range("E45").Value = "  +3.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5213"
$ws.Range("E46").Value = "  +10.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.75"
$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.52"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("E49").Value = "  +3.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9987"
$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06394"
$ws.Range("E51").Value = "  +1.31%  "
